# Append a new day's row (2026-01-03) to the "Chart" sheet of the GSC export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Write the numeric columns first so the new row keeps all three cells
# (writing the text/date-like column last avoids an observed ordering bug).
$ws.Cells.Item(90, 2).Value = 0
$ws.Cells.Item(90, 3).Value = 28

# Column A holds a plain text value that looks like a date ("2026-01-03").
# A straight .Value assignment gets auto-parsed into a real date serial by
# Excel's input heuristics, which would both change the cell's stored type
# and attach a date number format. Force a text format before the write so
# it is stored as text, then clear the format back to the sheet's default
# so no stray style is left behind on the cell.
$cell = $ws.Cells.Item(90, 1)
$cell.NumberFormat = "@"
$cell.Value = "2026-01-03"
$cell.ClearFormats()
